$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (run_time) new values for rows 2-11
$cVals = New-Object 'object[,]' 10,1
$cVals[0,0] = 0.5379502773284912
$cVals[1,0] = 0.5468440055847168
$cVals[2,0] = 0.6373727321624756
$cVals[3,0] = 0.5726428031921387
$cVals[4,0] = 0.5468497276306152
$cVals[5,0] = 0.5625274181365967
$cVals[6,0] = 0.5312473773956299
$cVals[7,0] = 0.5312507152557373
$cVals[8,0] = 0.5468757152557373
$cVals[9,0] = 0.5468735694885254
$ws.Range("C2:C11").Value = $cVals

# Columns E:Y (max_er + iter 0..19) new values for rows 2-11
$eyVals = New-Object 'object[,]' 10,21
$eyVals[0,0] = 1362.879238935046
$eyVals[0,1] = 0.08763877066616421
$eyVals[0,2] = 0.0659686939257412
$eyVals[0,3] = 0.05415184939863833
$eyVals[0,4] = 0.04777356964298513
$eyVals[0,5] = 0.04209803594138108
$eyVals[0,6] = 0.03931961072044678
$eyVals[0,7] = 0.036335348939463
$eyVals[0,8] = 0.03371375395731176
$eyVals[0,9] = 0.03169166742758392
$eyVals[0,10] = 0.03092998878537609
$eyVals[0,11] = 0.02872268616742391
$eyVals[0,12] = 0.02868768712854163
$eyVals[0,13] = 0.02804791659761854
$eyVals[0,14] = 0.02775839900052697
$eyVals[0,15] = 0.02742037517732443
$eyVals[0,16] = 0.02713258567321796
$eyVals[0,17] = 0.02688026967745171
$eyVals[0,18] = 0.02670058892456267
$eyVals[0,19] = 0.02665466711034592
$eyVals[0,20] = 0.02656684676286638
$eyVals[1,0] = 1522.855548330979
$eyVals[1,1] = 0.08361686332864879
$eyVals[1,2] = 0.06423223509984911
$eyVals[1,3] = 0.05331887230971315
$eyVals[1,4] = 0.04778818701787041
$eyVals[1,5] = 0.04346189893065654
$eyVals[1,6] = 0.04071720165153045
$eyVals[1,7] = 0.03819074901077502
$eyVals[1,8] = 0.03531036696737656
$eyVals[1,9] = 0.03462112202391379
$eyVals[1,10] = 0.03299750715313022
$eyVals[1,11] = 0.03210460208562309
$eyVals[1,12] = 0.03147271764788988
$eyVals[1,13] = 0.03104447805117199
$eyVals[1,14] = 0.03083306409177649
$eyVals[1,15] = 0.03042183779829276
$eyVals[1,16] = 0.03023575642590364
$eyVals[1,17] = 0.02994225360485253
$eyVals[1,18] = 0.02984155714810646
$eyVals[1,19] = 0.02973538314249876
$eyVals[1,20] = 0.02968529333978516
$eyVals[2,0] = 1424.804695031242
$eyVals[2,1] = 0.08640468316215603
$eyVals[2,2] = 0.06733635468590429
$eyVals[2,3] = 0.0565765207963202
$eyVals[2,4] = 0.04773599730404315
$eyVals[2,5] = 0.04110232628151335
$eyVals[2,6] = 0.0386170366768795
$eyVals[2,7] = 0.03540348277681836
$eyVals[2,8] = 0.03409179806983577
$eyVals[2,9] = 0.03252989001615429
$eyVals[2,10] = 0.03132877753908755
$eyVals[2,11] = 0.03053147885008178
$eyVals[2,12] = 0.029798049838418
$eyVals[2,13] = 0.02942136100226323
$eyVals[2,14] = 0.02903172112738912
$eyVals[2,15] = 0.02862035906896673
$eyVals[2,16] = 0.02841107059161844
$eyVals[2,17] = 0.02826432080213943
$eyVals[2,18] = 0.02797031877465897
$eyVals[2,19] = 0.02782767106898802
$eyVals[2,20] = 0.02777397066337703
$eyVals[3,0] = 1334.61952884516
$eyVals[3,1] = 0.08466164772061893
$eyVals[3,2] = 0.0653563507779675
$eyVals[3,3] = 0.05189310557398728
$eyVals[3,4] = 0.0420133442655346
$eyVals[3,5] = 0.03926378465200384
$eyVals[3,6] = 0.03611035277400767
$eyVals[3,7] = 0.03338242622980554
$eyVals[3,8] = 0.03207948419195458
$eyVals[3,9] = 0.03062691710373522
$eyVals[3,10] = 0.02890188730833924
$eyVals[3,11] = 0.02823492821027294
$eyVals[3,12] = 0.02755394327352632
$eyVals[3,13] = 0.02724653630065522
$eyVals[3,14] = 0.02699015815157269
$eyVals[3,15] = 0.02670298475785676
$eyVals[3,16] = 0.02656354568836831
$eyVals[3,17] = 0.02619010980841111
$eyVals[3,18] = 0.0261600494083787
$eyVals[3,19] = 0.02609478420039806
$eyVals[3,20] = 0.02601597522115321
$eyVals[4,0] = 1362.695779015841
$eyVals[4,1] = 0.08148141783192024
$eyVals[4,2] = 0.0648621535885256
$eyVals[4,3] = 0.05050776190242444
$eyVals[4,4] = 0.03938509009910161
$eyVals[4,5] = 0.03891252715628523
$eyVals[4,6] = 0.03578975799524509
$eyVals[4,7] = 0.0340451984605006
$eyVals[4,8] = 0.03267302015811795
$eyVals[4,9] = 0.03155870263938178
$eyVals[4,10] = 0.03074643779051179
$eyVals[4,11] = 0.02975779914869726
$eyVals[4,12] = 0.02917031985495603
$eyVals[4,13] = 0.02842532408198684
$eyVals[4,14] = 0.02828845370361033
$eyVals[4,15] = 0.02783671111642836
$eyVals[4,16] = 0.02756630559113368
$eyVals[4,17] = 0.02721721033649855
$eyVals[4,18] = 0.02690803394605336
$eyVals[4,19] = 0.02667212698057454
$eyVals[4,20] = 0.02656327054611776
$eyVals[5,0] = 1427.536268220741
$eyVals[5,1] = 0.08137416062519151
$eyVals[5,2] = 0.06206261508413942
$eyVals[5,3] = 0.05166531693420968
$eyVals[5,4] = 0.0457009525546058
$eyVals[5,5] = 0.04195284247268454
$eyVals[5,6] = 0.03826442422591787
$eyVals[5,7] = 0.03583366798151108
$eyVals[5,8] = 0.03419606376162903
$eyVals[5,9] = 0.03172141654112618
$eyVals[5,10] = 0.03152267879000625
$eyVals[5,11] = 0.0304989982389703
$eyVals[5,12] = 0.02979901203876512
$eyVals[5,13] = 0.02930614093560772
$eyVals[5,14] = 0.02873972457558988
$eyVals[5,15] = 0.02845220114562827
$eyVals[5,16] = 0.02822787381400904
$eyVals[5,17] = 0.02812344629002854
$eyVals[5,18] = 0.02792208928011815
$eyVals[5,19] = 0.02785596244070385
$eyVals[5,20] = 0.02782721770410802
$eyVals[6,0] = 1369.49363109059
$eyVals[6,1] = 0.08366840322310459
$eyVals[6,2] = 0.06530783125672728
$eyVals[6,3] = 0.05175219412386486
$eyVals[6,4] = 0.04322901155093072
$eyVals[6,5] = 0.0405033377045876
$eyVals[6,6] = 0.03696016912160736
$eyVals[6,7] = 0.03455058336522786
$eyVals[6,8] = 0.03257386970215584
$eyVals[6,9] = 0.03096274792744835
$eyVals[6,10] = 0.03048850592079426
$eyVals[6,11] = 0.02947498738543278
$eyVals[6,12] = 0.02889301128722907
$eyVals[6,13] = 0.0282789946581179
$eyVals[6,14] = 0.02787434459961883
$eyVals[6,15] = 0.02764339016123757
$eyVals[6,16] = 0.02727148585279834
$eyVals[6,17] = 0.02704267510064059
$eyVals[6,18] = 0.02686006151935953
$eyVals[6,19] = 0.0267554262381233
$eyVals[6,20] = 0.02669578228246764
$eyVals[7,0] = 1343.889752077903
$eyVals[7,1] = 0.0866742594030602
$eyVals[7,2] = 0.06605318781023191
$eyVals[7,3] = 0.05270522917609941
$eyVals[7,4] = 0.04331936367573792
$eyVals[7,5] = 0.03847274768775324
$eyVals[7,6] = 0.03408585404114008
$eyVals[7,7] = 0.03275891905176725
$eyVals[7,8] = 0.03106931380365048
$eyVals[7,9] = 0.02980120135938038
$eyVals[7,10] = 0.02847970326480066
$eyVals[7,11] = 0.02764410616898378
$eyVals[7,12] = 0.02756512019493003
$eyVals[7,13] = 0.02704746752264118
$eyVals[7,14] = 0.02669376419959543
$eyVals[7,15] = 0.02651191051881675
$eyVals[7,16] = 0.02651191051881675
$eyVals[7,17] = 0.02644560075285921
$eyVals[7,18] = 0.02635080824615261
$eyVals[7,19] = 0.02633567656213838
$eyVals[7,20] = 0.02619668132705463
$eyVals[8,0] = 1360.035236524651
$eyVals[8,1] = 0.08175339939961866
$eyVals[8,2] = 0.05939275389805506
$eyVals[8,3] = 0.05216592886107629
$eyVals[8,4] = 0.04404828603862633
$eyVals[8,5] = 0.04031534320730675
$eyVals[8,6] = 0.03645246841428675
$eyVals[8,7] = 0.03380732865763134
$eyVals[8,8] = 0.03154331307519034
$eyVals[8,9] = 0.02960283905516569
$eyVals[8,10] = 0.02888392030516455
$eyVals[8,11] = 0.0287912123092457
$eyVals[8,12] = 0.02833821593620636
$eyVals[8,13] = 0.02755637453761902
$eyVals[8,14] = 0.02733473727608399
$eyVals[8,15] = 0.02710533514592152
$eyVals[8,16] = 0.02703616932830748
$eyVals[8,17] = 0.02678071210946995
$eyVals[8,18] = 0.02669910760695038
$eyVals[8,19] = 0.02657388465567731
$eyVals[8,20] = 0.0265114081193889
$eyVals[9,0] = 1458.047476610283
$eyVals[9,1] = 0.08099376482044855
$eyVals[9,2] = 0.06349216034661383
$eyVals[9,3] = 0.05536012407293445
$eyVals[9,4] = 0.04544874253695703
$eyVals[9,5] = 0.04269979952529807
$eyVals[9,6] = 0.03861602624336143
$eyVals[9,7] = 0.03628460305708918
$eyVals[9,8] = 0.03375505558600902
$eyVals[9,9] = 0.03207306961909756
$eyVals[9,10] = 0.03144171816093407
$eyVals[9,11] = 0.03073216207147363
$eyVals[9,12] = 0.03023585081120483
$eyVals[9,13] = 0.02976004404791795
$eyVals[9,14] = 0.02930660986304776
$eyVals[9,15] = 0.02911696927665261
$eyVals[9,16] = 0.0289195182716698
$eyVals[9,17] = 0.02870807738541157
$eyVals[9,18] = 0.0286143382852668
$eyVals[9,19] = 0.02849398262840424
$eyVals[9,20] = 0.02842197810156497
$ws.Range("E2:Y11").Value = $eyVals

Write-Output "Updated log write mode values for rows 2-11"